$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new mail-log entry as row 45 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(45, 1).Value = "Afmelding nieuwsbrief"
$logs.Cells.Item(45, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(45, 3).Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Cells.Item(45, 4).Value = "Afmelding"
$logs.Cells.Item(45, 6).Value = "2025-06-17 22:06:14"
$logs.Cells.Item(45, 7).Value = "Nee"

# Extend the existing conditional-formatting blocks to cover the new row,
# keeping all rules/dxf mappings untouched.
$logs.Range("D2:D44").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D45"))
$logs.Range("G2:G44").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G45"))

# --- "Dashboard" sheet: bump the "Afmelding" tally from 6 to 7 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(4, 2).Value = 7
